$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (R) to the right of the existing "2020" column (Q),
# carrying over the same formatting as the corresponding Q cells.
$ws.Range("Q4:Q8").Copy() | Out-Null
$ws.Range("R4:R8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(4, 18).Value = 2021
$ws.Cells.Item(5, 18).Value = 47.8
$ws.Cells.Item(6, 18).Value = 20.7
$ws.Cells.Item(7, 18).Value = 9.8
$ws.Cells.Item(8, 18).Value = 17.3

# Move the active selection to P10, matching the post-edit workbook state.
$ws.Range("P10").Select() | Out-Null
